$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Narrow column B (closest achievable stored width to the authored 11.28515625)
$ws.Columns.Item(2).ColumnWidth = 10.5

# Header row (row 1) A1:C1 -> centered (style moves from bold-left to bold-centered)
$ws.Range("A1:C1").HorizontalAlignment = -4108  # xlCenter

# Data rows (2-32) A:C -> centered
$ws.Range("A2:C32").HorizontalAlignment = -4108  # xlCenter

# Update the selection to match the authored state
$ws.Range("G22:J22").Select()
